$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell without Excel's
# automatic number coercion (and without touching cell styles).
# We build the literal as a text formula in a scratch cell, copy it,
# and paste-special (values only) into the destination -- this keeps
# the destination cell's style (s=) untouched and yields a pure string
# cell, matching how the source data was authored.
$scratch = $ws.Range("ZZ1")
function Set-TextValue {
    param($cellRef, $val)
    $escaped = $val -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" '69.046.48'
Set-TextValue "E2" '  -1.20%  '
Set-TextValue "D3" '3.516.40'
Set-TextValue "E3" '  -2.24%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  +0.03%  '
Set-TextValue "D5" '586.45'
Set-TextValue "E5" '  +1.79%  '
Set-TextValue "D6" '171.96'
Set-TextValue "E6" '  -1.26%  '
Set-TextValue "D7" '0.611'
Set-TextValue "E7" '  +1.11%  '
Set-TextValue "D8" '3.511.24'
Set-TextValue "E8" '  -1.86%  '
Set-TextValue "E9" '  -0.08%  '
Set-TextValue "D10" '0.189'
Set-TextValue "E10" '  -2.52%  '
Set-TextValue "D11" '6.83'
Set-TextValue "E11" '  +2.95%  '
Set-TextValue "D12" '0.581'
Set-TextValue "E12" '  -2.67%  '
Set-TextValue "D13" '47.34'
Set-TextValue "E13" '  -1.06%  '
Set-TextValue "D14" '0.0000276'
Set-TextValue "E14" '  -1.75%  '
Set-TextValue "D15" '4.077.71'
Set-TextValue "E15" '  -2.47%  '
Set-TextValue "D16" '8.50'
Set-TextValue "E16" '  -3.01%  '
Set-TextValue "D17" '628.13'
Set-TextValue "E17" '  -4.92%  '
Set-TextValue "D18" '69.055.46'
Set-TextValue "E18" '  -1.20%  '
Set-TextValue "D19" '3.514.89'
Set-TextValue "E19" '  -3.03%  '
Set-TextValue "D20" '0.123'
Set-TextValue "E20" '  +1.30%  '
Set-TextValue "D21" '17.40'
Set-TextValue "E21" '  -0.83%  '
Set-TextValue "D22" '11.13'
Set-TextValue "E22" '  -0.82%  '
Set-TextValue "D23" '0.887'
Set-TextValue "E23" '  -3.32%  '
Set-TextValue "D24" '15.93'
Set-TextValue "E24" '  -5.78%  '
Set-TextValue "D25" '97.08'
Set-TextValue "E25" '  -1.93%  '
Set-TextValue "D26" '3.81'
Set-TextValue "E26" '  -1.56%  '
Set-TextValue "E27" '  +0.26%  '
Set-TextValue "D28" '2.63'
Set-TextValue "E28" '  -3.96%  '
Set-TextValue "D29" '9.28'
Set-TextValue "E29" '  -5.74%  '
Set-TextValue "D30" '32.71'
Set-TextValue "E30" '  -4.32%  '
Set-TextValue "D31" '8.54'
Set-TextValue "E31" '  -3.40%  '
Set-TextValue "D32" '3.14'
Set-TextValue "E32" '  -4.67%  '
Set-TextValue "E33" '  -2.09%  '
Set-TextValue "D34" '6.95'
Set-TextValue "E34" '  -3.85%  '
Set-TextValue "D35" '636.63'
Set-TextValue "E35" '  +11.29%  '
Set-TextValue "D36" '10.77'
Set-TextValue "E36" '  -1.38%  '
Set-TextValue "E37" '  -9.66%  '
Set-TextValue "E38" '  -2.77%  '
Set-TextValue "D39" '57.26'
Set-TextValue "E39" '  -0.90%  '
Set-TextValue "E40" '  -0.36%  '
Set-TextValue "D41" '0.0454'
Set-TextValue "E41" '  +1.33%  '
Set-TextValue "D42" '0.136'
Set-TextValue "E42" '  -2.35%  '
Set-TextValue "D43" '3.382.13'
Set-TextValue "E43" '  -4.84%  '
Set-TextValue "D44" '0.328'
Set-TextValue "E44" '  -3.20%  '
Set-TextValue "B45" 'InjectiveProtocol'
Set-TextValue "C45" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D45" '32.84'
Set-TextValue "E45" '  -4.27%  '
Set-TextValue "B46" 'PEPE'
Set-TextValue "C46" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D46" '0.0₃0700'
Set-TextValue "E46" '  -3.85%  '
Set-TextValue "D47" '2.54'
Set-TextValue "E47" '  -4.31%  '
Set-TextValue "D48" '2.73'
Set-TextValue "E48" '  -5.11%  '
Set-TextValue "E49" '  -1.18%  '
Set-TextValue "D50" '132.16'
Set-TextValue "E50" '  -1.97%  '
Set-TextValue "E51" '  +14.29%  '

$scratch.ClearContents()
$excel.CutCopyMode = 0

